$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "H" column (third comparison run) values; column I holds
# the shared formula =H/B-1 and will recalculate automatically.
$ws.Range("H3").Value = 7750
$ws.Range("H4").Value = 908529
$ws.Range("H5").Value = 54092
$ws.Range("H7").Value = 106332
$ws.Range("H8").Value = 1670547
$ws.Range("H9").Value = 1395981
$ws.Range("H10").Value = 1328514
$ws.Range("H11").Value = 1058054
$ws.Range("H12").Value = 1453187
$ws.Range("H14").Value = 164137

# Match the saved selection/cursor position recorded in the workbook.
$ws.Range("I16").Select()
